$d = $word.ActiveDocument

# The literal "www.drpaulduenas.com" text in the footer becomes a dynamic
# MERGEFIELD ("=website") field result, matching the pattern already used
# by the neighboring "emergency_number" / "consultation.branch_office.*"
# merge fields in the same footer:
#
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="begin"/></w:r>
#   <w:r><w:rPr>...</w:rPr><w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="separate"/></w:r>
#   <w:r><w:rPr>...</w:rPr><w:t>�=website�</w:t></w:r>
#   <w:r><w:rPr>...</w:rPr><w:fldChar w:fldCharType="end"/></w:r>

$target = "www.drpaulduenas.com"

# Run-level formatting (w:rPr) carried by the original "www.drpaulduenas.com"
# run, and reused verbatim for every run of the new field.
$rPrXml = '<w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

# The paragraph that only contains "www.drpaulduenas.com" (identified by its
# w14:paraId in the source document) - its <w:pPr> / attributes are left
# untouched by the edit, only its run content changes.
$paraId = "24EA949D"
$paraAttrs = 'w14:paraId="24EA949D" w14:textId="77777777" w:rsidR="004D2A29" w:rsidRDefault="004D2A29" w:rsidP="004D2A29"'
$pPrXml = '<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Avenir Book" w:hAnsi="Avenir Book"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>'

$newParaXml = (
  '<w:p ' + $paraAttrs + '>' +
    $pPrXml +
    '<w:r>' + $rPrXml + '<w:fldChar w:fldCharType="begin"/></w:r>' +
    '<w:r>' + $rPrXml + '<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>' +
    '<w:r>' + $rPrXml + '<w:fldChar w:fldCharType="separate"/></w:r>' +
    '<w:r>' + $rPrXml + '<w:t>' + [char]0x00AB + '=website' + [char]0x00BB + '</w:t></w:r>' +
    '<w:r>' + $rPrXml + '<w:fldChar w:fldCharType="end"/></w:r>' +
  '</w:p>'
)

$fieldXml = (
  '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' + $newParaXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
)

# Fallback replacement XML used if the matched paragraph's identity can't be
# confirmed (keeps only the run-level formatting, letting Word regenerate
# a plain paragraph around it rather than risk writing wrong attributes).
$fallbackParaXml = (
  '<w:p>' +
    '<w:pPr><w:pStyle w:val="Footer"/></w:pPr>' +
    '<w:r>' + $rPrXml + '<w:fldChar w:fldCharType="begin"/></w:r>' +
    '<w:r>' + $rPrXml + '<w:instrText xml:space="preserve"> MERGEFIELD =website \* MERGEFORMAT </w:instrText></w:r>' +
    '<w:r>' + $rPrXml + '<w:fldChar w:fldCharType="separate"/></w:r>' +
    '<w:r>' + $rPrXml + '<w:t>' + [char]0x00AB + '=website' + [char]0x00BB + '</w:t></w:r>' +
    '<w:r>' + $rPrXml + '<w:fldChar w:fldCharType="end"/></w:r>' +
  '</w:p>'
)
$fallbackFieldXml = (
  '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $fallbackParaXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
)

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers.Item($i)
        if (-not $ftr.Exists) { continue }

        $searchRange = $ftr.Range.Duplicate
        while ($searchRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)) {
            $matchRange = $searchRange.Duplicate

            # Confirm this run's paragraph is the one we expect (full text of
            # the paragraph is exactly the target URL) before rewriting the
            # whole paragraph with the known-good attributes; otherwise fall
            # back to a generic wrapper so we never destroy unrelated content.
            $para = $matchRange.Paragraphs.Item(1)
            $useKnownParagraph = ($para.Range.Text.Trim() -eq $target)

            # Clear the matched text; this collapses matchRange to a single
            # insertion point in place of the old run.
            $matchRange.Text = ""

            if ($useKnownParagraph) {
                $matchRange.InsertXML($fieldXml)
            } else {
                $matchRange.InsertXML($fallbackFieldXml)
            }

            # Resume searching right after the footer content we just wrote.
            $searchRange = $ftr.Range.Duplicate
            $searchRange.Start = $ftr.Range.End
        }
    }
}
